$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update unit abbreviation "kVArh" -> "kvarh" (all cells sharing this value)
$ws.Range("C11").Value = "kvarh"
$ws.Range("C12").Value = "kvarh"
$ws.Range("C13").Value = "kvarh"
$ws.Range("C14").Value = "kvarh"
$ws.Range("C15").Value = "kvarh"
$ws.Range("C39").Value = "kvarh"
$ws.Range("C40").Value = "kvarh"
$ws.Range("C41").Value = "kvarh"
$ws.Range("C42").Value = "kvarh"

# Update matching descriptions "kVAR-hours" -> "kvar-hours"
$ws.Range("D11").Value = "kvar-hours, Negative - Reactive Energy capacitive exported"
$ws.Range("D12").Value = "kvar-hours, Positive - Reactive Energy capacitive imported"
$ws.Range("D13").Value = "kvar-hours, Positive - Reactive Energy inductive exported"
$ws.Range("D14").Value = "kvar-hours, Positive - Reactive Energy inductive imported"
$ws.Range("D15").Value = "kvar-hours - Reactive Energy total"

# Update unit "VAr" -> "var"
$ws.Range("C28").Value = "var"

Write-Output "edit complete"
